# 2017 Fix to date in powerpoint
#
# The title on slide 1 reads "FME 2016 Training" and needs to become
# "FME 2017 Training". The author made the edit by retyping the year
# in place (selecting just the "2016" portion and typing "2017"),
# which is why the authored OOXML ends up with the run split into
# three separate <a:r> runs ("FME ", "2017 ", "Training") that all
# share identical run properties, rather than one single run.
#
# We reproduce that exact run split here by writing the full new
# text first and then re-touching each of the three sub-ranges (via
# TextRange.Characters(start, length)) so PowerPoint keeps them as
# distinct runs in the saved XML, instead of collapsing them back
# into a single run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

# Replace the whole title text first.
$tr.Text = "FME 2017 Training"

# Re-apply formatting per sub-range so the saved file keeps three
# separate runs, matching the authored edit.
$run1 = $tr.Characters(1, 4)          # "FME "
$run1.Font.Size = 60

$run2 = $tr.Characters(5, 5)          # "2017 "
$run2.Font.Size = 60

$run3 = $tr.Characters(10, 8)         # "Training"
$run3.Font.Size = 60
